# Weekly fruit/vegetable price update: a new week's record is inserted at
# the top of the data table (row 21, right after the table's header block
# that ends at row 20), which pushes all the existing data rows (21-94)
# down by one (to 22-95). The workbook's dimension grows from A1:T94 to
# A1:T95 automatically as a result of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 21 - shifts rows 21:94 down to 22:95.
$ws.Rows("21").Insert()

# Populate the newly inserted row 21 with this week's record.
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = 44600
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100108
$ws.Range("H21").Value = "Tropicales y subtropicales"
$ws.Range("I21").Value = 100108002
$ws.Range("J21").Value = "Mango"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 7000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 7500
$ws.Range("Q21").Value = "`$/bandeja 4 kilos"
$ws.Range("R21").Value = "Perú"
$ws.Range("S21").Value = 1875
$ws.Range("T21").Value = 4
